# Weekly refresh: two new sampling rows (2021-11-03) are inserted right
# after the first 5 rows of the sheet (i.e. become the new rows 7-8),
# pushing every following row down by two positions. The two oldest rows
# that fall off the bottom (old rows 113-114) are appended as the new
# rows 115-116, so no data is lost - it's simply re-paginated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$old = $used.Value2

$oldRows = $old.GetUpperBound(0)   # 114
$cols = $old.GetUpperBound(1)      # 18
$newRows = $oldRows + 2            # 116

# New data for the two freshly-sampled rows (old "row 7" slot becomes
# row 9, etc. - see below).
$row7 = @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44503, 13, 100112001, "Berenjena", "Sin especificar", "Primera", 220, 7000, 8000, 7545, "`$/caja 50 unidades", "Región de Arica y Parinacota", 151, 50, "Hortaliza")
$row8 = @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44503, 13, 100112001, "Berenjena", "Sin especificar", "Primera", 90, 12000, 13000, 12444, "`$/caja 70 unidades", "Provincia de Huasco", 178, 70, "Hortaliza")

$new = New-Object 'object[,]' $newRows, $cols

for ($c = 1; $c -le $cols; $c++) {
    # Rows 1-6 (header + first 5 data rows) are untouched.
    for ($r = 1; $r -le 6; $r++) {
        $new[$r - 1, $c - 1] = $old[$r, $c]
    }

    # New rows 7 and 8.
    $new[6, $c - 1] = $row7[$c - 1]
    $new[7, $c - 1] = $row8[$c - 1]

    # Old rows 7..114 slide down two rows to become new rows 9..116.
    for ($r = 7; $r -le $oldRows; $r++) {
        $new[$r + 1, $c - 1] = $old[$r, $c]
    }
}

$target = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($newRows, $cols))
$target.Value = $new

# The two brand-new rows at the bottom (116th/115th... i.e. the shifted-
# down tail) need the date number format re-applied on column D since
# they are genuinely new cells with no prior formatting.
$ws.Range("D115:D116").NumberFormat = "YYYY-MM-DD HH:MM:SS"
